function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Cells.Item(2, 4) "37.428.90"
Set-TextValue $ws.Cells.Item(2, 5) "  +2.50%  "
Set-TextValue $ws.Cells.Item(3, 4) "2.099.92"
Set-TextValue $ws.Cells.Item(3, 5) "  +4.51%  "
Set-TextValue $ws.Cells.Item(4, 5) "  +0.03%  "
Set-TextValue $ws.Cells.Item(5, 4) "250.95"
Set-TextValue $ws.Cells.Item(5, 5) "  +2.15%  "
Set-TextValue $ws.Cells.Item(6, 4) "0.663"
Set-TextValue $ws.Cells.Item(6, 5) "  +1.30%  "
Set-TextValue $ws.Cells.Item(7, 5) "  -0.07%  "
Set-TextValue $ws.Cells.Item(8, 4) "51.14"
Set-TextValue $ws.Cells.Item(8, 5) "  +14.06%  "
Set-TextValue $ws.Cells.Item(9, 4) "61.55"
Set-TextValue $ws.Cells.Item(9, 5) "  +10.50%  "
Set-TextValue $ws.Cells.Item(10, 4) "0.374"
Set-TextValue $ws.Cells.Item(10, 5) "  +3.85%  "
Set-TextValue $ws.Cells.Item(11, 4) "0.0745"
Set-TextValue $ws.Cells.Item(11, 5) "  +3.89%  "
Set-TextValue $ws.Cells.Item(12, 5) "  +6.39%  "
Set-TextValue $ws.Cells.Item(13, 4) "15.39"
Set-TextValue $ws.Cells.Item(13, 5) "  +7.37%  "
Set-TextValue $ws.Cells.Item(14, 4) "2.403.54"
Set-TextValue $ws.Cells.Item(14, 5) "  +4.36%  "
Set-TextValue $ws.Cells.Item(15, 4) "0.833"
Set-TextValue $ws.Cells.Item(15, 5) "  +4.50%  "
Set-TextValue $ws.Cells.Item(16, 4) "2.101.58"
Set-TextValue $ws.Cells.Item(16, 5) "  +4.48%  "
Set-TextValue $ws.Cells.Item(17, 4) "5.12"
Set-TextValue $ws.Cells.Item(17, 5) "  +5.02%  "
Set-TextValue $ws.Cells.Item(18, 4) "37.293.10"
Set-TextValue $ws.Cells.Item(18, 5) "  +1.90%  "
Set-TextValue $ws.Cells.Item(19, 4) "72.20"
Set-TextValue $ws.Cells.Item(19, 5) "  +1.86%  "
Set-TextValue $ws.Cells.Item(20, 4) "0.0₃0832"
Set-TextValue $ws.Cells.Item(20, 5) "  +2.46%  "
Set-TextValue $ws.Cells.Item(21, 4) "13.57"
Set-TextValue $ws.Cells.Item(21, 5) "  +5.11%  "
Set-TextValue $ws.Cells.Item(22, 4) "240.38"
Set-TextValue $ws.Cells.Item(22, 5) "  +2.70%  "
Set-TextValue $ws.Cells.Item(23, 5) "  +5.01%  "
Set-TextValue $ws.Cells.Item(24, 5) "  +0.18%  "
Set-TextValue $ws.Cells.Item(25, 4) "2.45"
Set-TextValue $ws.Cells.Item(25, 5) "  +1.48%  "
Set-TextValue $ws.Cells.Item(26, 4) "170.09"
Set-TextValue $ws.Cells.Item(26, 5) "  +5.34%  "
Set-TextValue $ws.Cells.Item(27, 4) "9.17"
Set-TextValue $ws.Cells.Item(27, 5) "  +8.38%  "
Set-TextValue $ws.Cells.Item(28, 4) "20.66"
Set-TextValue $ws.Cells.Item(28, 5) "  +5.73%  "
Set-TextValue $ws.Cells.Item(29, 4) "2.00"
Set-TextValue $ws.Cells.Item(29, 5) "  +1.06%  "
Set-TextValue $ws.Cells.Item(30, 4) "0.123"
Set-TextValue $ws.Cells.Item(30, 5) "  +0.63%  "
Set-TextValue $ws.Cells.Item(31, 4) "1.06"
Set-TextValue $ws.Cells.Item(31, 5) "  +25.19%  "
Set-TextValue $ws.Cells.Item(32, 4) "4.49"
Set-TextValue $ws.Cells.Item(32, 5) "  +3.53%  "
Set-TextValue $ws.Cells.Item(33, 4) "0.0608"
Set-TextValue $ws.Cells.Item(33, 5) "  +4.87%  "
Set-TextValue $ws.Cells.Item(34, 5) "  +13.88%  "
Set-TextValue $ws.Cells.Item(35, 5) "  -0.10%  "
Set-TextValue $ws.Cells.Item(36, 4) "2.33"
Set-TextValue $ws.Cells.Item(36, 5) "  +10.63%  "
Set-TextValue $ws.Cells.Item(37, 4) "19.20"
Set-TextValue $ws.Cells.Item(37, 5) "  -6.74%  "
Set-TextValue $ws.Cells.Item(38, 4) "4.09"
Set-TextValue $ws.Cells.Item(38, 5) "  +1.73%  "
Set-TextValue $ws.Cells.Item(39, 4) "1.83"
Set-TextValue $ws.Cells.Item(39, 5) "  -1.15%  "
Set-TextValue $ws.Cells.Item(40, 4) "1.32"
Set-TextValue $ws.Cells.Item(40, 5) "  -0.79%  "
Set-TextValue $ws.Cells.Item(41, 4) "17.97"
Set-TextValue $ws.Cells.Item(41, 5) "  +11.31%  "
Set-TextValue $ws.Cells.Item(42, 4) "0.0224"
Set-TextValue $ws.Cells.Item(42, 5) "  +4.27%  "
Set-TextValue $ws.Cells.Item(43, 5) "  +8.83%  "
Set-TextValue $ws.Cells.Item(44, 4) "98.98"
Set-TextValue $ws.Cells.Item(44, 5) "  +2.91%  "
Set-TextValue $ws.Cells.Item(45, 2) "Cronos"
Set-TextValue $ws.Cells.Item(45, 3) "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Cells.Item(45, 4) "0.0899"
Set-TextValue $ws.Cells.Item(45, 5) "  +10.68%  "
Set-TextValue $ws.Cells.Item(46, 2) "HuobiToken"
Set-TextValue $ws.Cells.Item(46, 3) "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Cells.Item(46, 4) "2.73"
Set-TextValue $ws.Cells.Item(46, 5) "  +0.29%  "
Set-TextValue $ws.Cells.Item(47, 5) "  +10.03%  "
Set-TextValue $ws.Cells.Item(48, 4) "1.321.48"
Set-TextValue $ws.Cells.Item(48, 5) "  +1.41%  "
Set-TextValue $ws.Cells.Item(49, 4) "6.98"
Set-TextValue $ws.Cells.Item(49, 5) "  +14.63%  "
Set-TextValue $ws.Cells.Item(50, 4) "2.280.01"
Set-TextValue $ws.Cells.Item(50, 5) "  +4.20%  "
Set-TextValue $ws.Cells.Item(51, 4) "2.29"
Set-TextValue $ws.Cells.Item(51, 5) "  +3.31%  "
